# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
# (the Simulate_Season.py portion lives outside this workbook; here we only
# apply the Week 17 stat updates to the "Players Data.xlsx" sheets.)

$wb = $excel.ActiveWorkbook

# --- Rushing sheet -------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# J.Allen
$rushing.Range("C2").Value = 33
$rushing.Range("D2").Value = 35
$rushing.Range("E2").Value = 36
$rushing.Range("F2").Value = 30

# D.Singletary
$rushing.Range("C3").Value = 91
$rushing.Range("D3").Value = 68
$rushing.Range("F3").Value = 33

# Z.Moss
$rushing.Range("C4").Value = 45
$rushing.Range("D4").Value = 33

# I.McKenzie
$rushing.Range("D8").Value = 3
$rushing.Range("F8").Value = 4

# --- Receiving sheet -------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# D.Singletary
$receiving.Range("C2").Value = 45

# S.Diggs
$receiving.Range("C6").Value = 108
$receiving.Range("D6").Value = 77
$receiving.Range("E6").Value = 31
$receiving.Range("G6").Value = 26
$receiving.Range("H6").Value = 15

# C.Beasley
$receiving.Range("C8").Value = 98
$receiving.Range("D8").Value = 73
$receiving.Range("E8").Value = 10
$receiving.Range("G8").Value = 15

# G.Davis
$receiving.Range("C9").Value = 29
$receiving.Range("D9").Value = 19
$receiving.Range("G9").Value = 14
$receiving.Range("H9").Value = 8

# I.McKenzie
$receiving.Range("C10").Value = 11
$receiving.Range("D10").Value = 8
$receiving.Range("G10").Value = 2
$receiving.Range("H10").Value = 2

# D.Knox
$receiving.Range("C12").Value = 46
$receiving.Range("E12").Value = 18

# The Rushing tab ends up the active/selected sheet after this edit.
$rushing.Activate()
